$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row is inserted at row 135 (pushing the
# existing rows 135..222 down to 136..223, growing the used range to
# A1:R223). Insert a fresh row there, shifting everything below it down.
$ws.Range("A135:R135").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A135").Value = 10
$ws.Range("B135").Value = "Vega Modelo de Temuco"
$ws.Range("C135").Value = "La Araucanía"
$ws.Range("D135").Value = 44777
$ws.Range("E135").Value = 9
$ws.Range("F135").Value = 100112005
$ws.Range("G135").Value = "Puerro"
$ws.Range("H135").Value = "Azul de Maquehue"
$ws.Range("I135").Value = "Primera"
$ws.Range("J135").Value = 85
$ws.Range("K135").Value = 16000
$ws.Range("L135").Value = 16000
$ws.Range("M135").Value = 16000
$ws.Range("N135").Value = "$/docena de paquetes"
$ws.Range("O135").Value = "Provincia de Cautín"
$ws.Range("P135").Value = 1333
$ws.Range("Q135").Value = 12
$ws.Range("R135").Value = "Hortaliza"
